$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Drop the last four report rows (old rows 9-12) - the report now
#    only covers the 25/28-FEB and 28-MAR runs (rows 2-8).
# ------------------------------------------------------------------
$ws.Rows("9:12").Delete()

# ------------------------------------------------------------------
# 2) Refresh the date column for the rows whose run-date moved on.
#    Force a text format first so Excel does not reinterpret the
#    "DD-MMM-YY" strings as real dates (the sheet stores them as
#    plain text, matching the original inline-string cells).
# ------------------------------------------------------------------
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("A2:A4").Value = "28-FEB-26"

$ws.Range("A6:A8").NumberFormat = "@"
$ws.Range("A6:A8").Value = "28-MAR-26"

# Re-apply the plain bordered look (style used by every other data
# cell) to the cells we just forced to text, since changing
# NumberFormat pushed them onto a new style slot.
$ws.Range("K2").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$ws.Range("A6:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3) Row 2 - Air Arabia Egypt E5-592 vs SM-436 (28-FEB-26 run)
# ------------------------------------------------------------------
$ws.Range("D2").Value = 579
$ws.Range("E2").Value = 1793
$ws.Range("F2").Value = -1214
$ws.Range("J2").Value = "HIGH THREAT ALERT - NEED ACTION"

# ------------------------------------------------------------------
# 4) Row 3 - Nesma Airlines NE-141 vs SM-436 (28-FEB-26 run)
# ------------------------------------------------------------------
$ws.Range("D3").Value = 600
$ws.Range("E3").Value = 1793
$ws.Range("F3").Value = -1193
$ws.Range("J3").Value = "HIGH THREAT ALERT - NEED ACTION"

# ------------------------------------------------------------------
# 5) Row 4 - Nile Air NP-116 vs SM-436 (28-FEB-26 run)
# ------------------------------------------------------------------
$ws.Range("C4").Value = "Nile Air NP-116"
$ws.Range("D4").Value = 600
$ws.Range("E4").Value = 1793
$ws.Range("F4").Value = -1193
$ws.Range("J4").Value = "HIGH THREAT ALERT - NEED ACTION"

# Rows 2-4 were "LOW THREAT" (green) before; they are now high-threat
# rows, so pick up the red high-threat look already used on row 5.
$ws.Range("J5").Copy()
$ws.Range("J2:J4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 6) Row 5 - EgyptAir MS-634 vs SM-436 (28-FEB-26 run)
# ------------------------------------------------------------------
$ws.Range("C5").Value = "EgyptAir MS-634"
$ws.Range("D5").Value = 1178
$ws.Range("F5").Value = -615
$ws.Range("G5").Value = 46
$ws.Range("I5").Value = -16

# ------------------------------------------------------------------
# 7) Row 6 - Air Arabia Egypt E5-592 vs SM-436 (28-MAR-26 run)
# ------------------------------------------------------------------
$ws.Range("C6").Value = "Air Arabia Egypt E5-592"
$ws.Range("D6").Value = 1020
$ws.Range("E6").Value = 2183
$ws.Range("F6").Value = -1163

# ------------------------------------------------------------------
# 8) Row 7 - Nile Air NP-116 vs SM-436 (28-MAR-26 run)
# ------------------------------------------------------------------
$ws.Range("C7").Value = "Nile Air NP-116"
$ws.Range("D7").Value = 1080
$ws.Range("E7").Value = 2183
$ws.Range("F7").Value = -1103

# ------------------------------------------------------------------
# 9) Row 8 - EgyptAir MS-634 vs SM-436 (28-MAR-26 run)
# ------------------------------------------------------------------
$ws.Range("D8").Value = 1495
$ws.Range("E8").Value = 2183
$ws.Range("F8").Value = -688
